# Update the "dSF" (column F) values for the crawford_kutter dataset sheet.
# This reflects a repull of data / recalculated mean, which changed several
# dSF values from their previous (stale) numbers to the freshly computed ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 1
    "F4"  = -4
    "F6"  = -5
    "F7"  = 1
    "F8"  = -1
    "F9"  = -2
    "F10" = 3
    "F12" = 2
    "F13" = -3
    "F16" = -8
    "F17" = -3
    "F20" = 4
    "F21" = 1
    "F25" = -3
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
